$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1820728291316527
$ws.Range("C2").Value = 0.5882352941176471
$ws.Range("J2").Value = 0.005602240896358543
$ws.Range("O2").Value = 0.002801120448179272
$ws.Range("P2").Value = 0.123249299719888
$ws.Range("S2").Value = 0.09803921568627451
$ws.Range("B3").Value = 0.009049773755656109
$ws.Range("C3").Value = 0.04072398190045249
$ws.Range("J3").Value = 0.02262443438914027
$ws.Range("P3").Value = 0.7420814479638009
$ws.Range("S3").Value = 0.1855203619909502
$ws.Range("J4").Value = 0.01923076923076923
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.2307692307692308
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.4
$ws.Range("B6").Value = 0.04460966542750929
$ws.Range("D6").Value = 0.01115241635687732
$ws.Range("F6").Value = 0.02973977695167286
$ws.Range("J6").Value = 0.2788104089219331
$ws.Range("O6").Value = 0.02230483271375465
$ws.Range("Q6").Value = 0.1672862453531599
$ws.Range("R6").Value = 0.07434944237918216
$ws.Range("S6").Value = 0.3717472118959108
$ws.Range("B7").Value = 0.1367924528301887
$ws.Range("D7").Value = 0.01886792452830189
$ws.Range("F7").Value = 0.04716981132075472
$ws.Range("J7").Value = 0.1556603773584906
$ws.Range("O7").Value = 0.02358490566037736
$ws.Range("Q7").Value = 0.1933962264150944
$ws.Range("R7").Value = 0.05660377358490566
$ws.Range("S7").Value = 0.3679245283018868
$ws.Range("B8").Value = 0.09797822706065319
$ws.Range("D8").Value = 0.02332814930015552
$ws.Range("E8").Value = 0.003110419906687403
$ws.Range("F8").Value = 0.07776049766718507
$ws.Range("J8").Value = 0.1259720062208398
$ws.Range("O8").Value = 0.02021772939346812
$ws.Range("Q8").Value = 0.1866251944012442
$ws.Range("R8").Value = 0.05909797822706065
$ws.Range("S8").Value = 0.4059097978227061
$ws.Range("B9").Value = 0.1149425287356322
$ws.Range("D9").Value = 0.005747126436781609
$ws.Range("E9").Value = 0.005747126436781609
$ws.Range("F9").Value = 0.08045977011494253
$ws.Range("J9").Value = 0.1264367816091954
$ws.Range("O9").Value = 0.01149425287356322
$ws.Range("Q9").Value = 0.2011494252873563
$ws.Range("R9").Value = 0.07471264367816093
$ws.Range("S9").Value = 0.3793103448275862
$ws.Range("B10").Value = 0.1278982797307405
$ws.Range("D10").Value = 0.02243829468960359
$ws.Range("E10").Value = 0.002243829468960359
$ws.Range("F10").Value = 0.07105459985041136
$ws.Range("J10").Value = 0.1151832460732984
$ws.Range("O10").Value = 0.01869857890800299
$ws.Range("Q10").Value = 0.2146596858638743
$ws.Range("R10").Value = 0.06357516828721017
$ws.Range("S10").Value = 0.3642483171278983
$ws.Range("G11").Value = 0.1316614420062696
$ws.Range("J11").Value = 0.1379310344827586
$ws.Range("K11").Value = 0.1849529780564263
$ws.Range("L11").Value = 0.5235109717868338
$ws.Range("S11").Value = 0.0219435736677116
$ws.Range("G12").Value = 0.7120418848167539
$ws.Range("J12").Value = 0.162303664921466
$ws.Range("K12").Value = 0.01570680628272251
$ws.Range("L12").Value = 0.03664921465968586
$ws.Range("S12").Value = 0.07329842931937172
$ws.Range("G13").Value = 0.7413793103448276
$ws.Range("J13").Value = 0.1206896551724138
$ws.Range("S13").Value = 0.1379310344827586
$ws.Range("F15").Value = 0.03658536585365853
$ws.Range("H15").Value = 0.2479674796747967
$ws.Range("I15").Value = 0.02845528455284553
$ws.Range("J15").Value = 0.2479674796747967
$ws.Range("K15").Value = 0.06097560975609756
$ws.Range("M15").Value = 0.008130081300813009
$ws.Range("O15").Value = 0.06504065040650407
$ws.Range("S15").Value = 0.3048780487804878
$ws.Range("F16").Value = 0.00816326530612245
$ws.Range("H16").Value = 0.1795918367346939
$ws.Range("I16").Value = 0.05714285714285714
$ws.Range("J16").Value = 0.363265306122449
$ws.Range("K16").Value = 0.1346938775510204
$ws.Range("M16").Value = 0.01224489795918367
$ws.Range("O16").Value = 0.06122448979591837
$ws.Range("S16").Value = 0.1836734693877551
$ws.Range("F17").Value = 0.01518026565464896
$ws.Range("H17").Value = 0.2580645161290323
$ws.Range("I17").Value = 0.07020872865275142
$ws.Range("J17").Value = 0.3681214421252372
$ws.Range("K17").Value = 0.1043643263757116
$ws.Range("M17").Value = 0.02656546489563567
$ws.Range("O17").Value = 0.05502846299810247
$ws.Range("S17").Value = 0.1024667931688805
$ws.Range("F18").Value = 0.03571428571428571
$ws.Range("H18").Value = 0.1666666666666667
$ws.Range("I18").Value = 0.05357142857142857
$ws.Range("J18").Value = 0.3869047619047619
$ws.Range("K18").Value = 0.07738095238095238
$ws.Range("M18").Value = 0.0119047619047619
$ws.Range("O18").Value = 0.08928571428571429
$ws.Range("S18").Value = 0.1785714285714286
$ws.Range("F19").Value = 0.02812071330589849
$ws.Range("H19").Value = 0.2578875171467764
$ws.Range("I19").Value = 0.06927297668038408
$ws.Range("J19").Value = 0.3429355281207133
$ws.Range("K19").Value = 0.09327846364883402
$ws.Range("M19").Value = 0.02400548696844993
$ws.Range("N19").Value = 0.001371742112482853
$ws.Range("O19").Value = 0.06515775034293553
$ws.Range("S19").Value = 0.1179698216735254

Write-Output "Updated 113 cells"
